# Updates the "Estado de Cuenta" (account statement) database:
#  - Column E (Periodo Mora) for rows 16-48 is reordered from descending
#    (2407 .. 2111) to ascending (2111 .. 2407) chronological order.
#  - Column F (Valor Mora) keeps the same 36360 values for all periods
#    except the new last period (now row 48) which keeps the 27876 value
#    that used to belong to the first period.
#  - Column G (Salario Basico) is updated to the new salary value 909000
#    for every period row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "2111", "2112",
    "2201", "2202", "2203", "2204", "2205", "2206", "2207", "2208", "2209", "2210", "2211", "2212",
    "2301", "2302", "2303", "2304", "2305", "2306", "2307", "2308", "2309", "2310", "2311", "2312",
    "2401", "2402", "2403", "2404", "2405", "2406", "2407"
)

$firstRow = 16
$lastRow = 48

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]   # column E - Periodo Mora
    $ws.Cells.Item($row, 7).Value = 909000          # column G - Salario Basico
    if ($row -eq $lastRow) {
        $ws.Cells.Item($row, 6).Value = 27876        # column F - Valor Mora (moved to last row)
    } else {
        $ws.Cells.Item($row, 6).Value = 36360        # column F - Valor Mora
    }
}
